$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells retain text formatting (matching original inlineStr values)
# then set the updated values from the latest crypto data refresh.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.037.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.916.04"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.39"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.06"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.94"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.97%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.42"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.398.40"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.934.77"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.915.98"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "433.39"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.39"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.85"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.32%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.97"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.95%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.122"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.51%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "41.82"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.287"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "375.78"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.690.96"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "132.98"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.80"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.124"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.52%  "
